# Adds a new "LEAD" worksheet (pivot-style summary of Lead_Source counts by
# year) after the existing "PRIOR" sheet, matching the author's commit:
# "added pivot tables for LEAD categories".

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)   # "PRIOR" - used as a formatting template

# New sheet goes at the end of the tab strip, becomes the active sheet/tab
# (mirrors Excel's own behaviour: Worksheets.Add() w/ After: = lastSheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "LEAD"

# ---- header row -----------------------------------------------------
$newSheet.Cells.Item(1,1).Value = "Lead_Source"
$newSheet.Cells.Item(1,2).Value = 2017
$newSheet.Cells.Item(1,3).Value = 2018
$newSheet.Cells.Item(1,4).Value = 2019

# ---- data rows (Lead_Source counts per year) -------------------------
$newSheet.Cells.Item(2,1).Value = "College Fair"
$newSheet.Cells.Item(2,2).Value = 7783
$newSheet.Cells.Item(2,3).Value = 5208
$newSheet.Cells.Item(2,4).Value = 5604
$newSheet.Cells.Item(3,1).Value = "College Visit"
$newSheet.Cells.Item(3,2).Value = 1
$newSheet.Cells.Item(3,3).Value = 12
$newSheet.Cells.Item(3,4).Value = 0
$newSheet.Cells.Item(4,1).Value = "Education Fair"
$newSheet.Cells.Item(4,2).Value = 3
$newSheet.Cells.Item(4,3).Value = 2
$newSheet.Cells.Item(4,4).Value = 10
$newSheet.Cells.Item(5,1).Value = "Email"
$newSheet.Cells.Item(5,2).Value = 116
$newSheet.Cells.Item(5,3).Value = 24
$newSheet.Cells.Item(5,4).Value = 8
$newSheet.Cells.Item(6,1).Value = "Event Registration"
$newSheet.Cells.Item(6,2).Value = 8074
$newSheet.Cells.Item(6,3).Value = 9131
$newSheet.Cells.Item(6,4).Value = 8803
$newSheet.Cells.Item(7,1).Value = "Group Tour"
$newSheet.Cells.Item(7,2).Value = 252
$newSheet.Cells.Item(7,3).Value = 279
$newSheet.Cells.Item(7,4).Value = 206
$newSheet.Cells.Item(8,1).Value = "High School Visit"
$newSheet.Cells.Item(8,2).Value = 5427
$newSheet.Cells.Item(8,3).Value = 5410
$newSheet.Cells.Item(8,4).Value = 42
$newSheet.Cells.Item(9,1).Value = "Initiated Application"
$newSheet.Cells.Item(9,2).Value = 9010
$newSheet.Cells.Item(9,3).Value = 8
$newSheet.Cells.Item(9,4).Value = 0
$newSheet.Cells.Item(10,1).Value = "International Travel"
$newSheet.Cells.Item(10,2).Value = 1
$newSheet.Cells.Item(10,3).Value = 0
$newSheet.Cells.Item(10,4).Value = 0
$newSheet.Cells.Item(11,1).Value = "Lead Card"
$newSheet.Cells.Item(11,2).Value = 148
$newSheet.Cells.Item(11,3).Value = 175
$newSheet.Cells.Item(11,4).Value = 145
$newSheet.Cells.Item(12,1).Value = "Office Visit"
$newSheet.Cells.Item(12,2).Value = 200
$newSheet.Cells.Item(12,3).Value = 131
$newSheet.Cells.Item(12,4).Value = 240
$newSheet.Cells.Item(13,1).Value = "Other"
$newSheet.Cells.Item(13,2).Value = 1208
$newSheet.Cells.Item(13,3).Value = 352
$newSheet.Cells.Item(13,4).Value = 1
$newSheet.Cells.Item(14,1).Value = "Phone"
$newSheet.Cells.Item(14,2).Value = 17
$newSheet.Cells.Item(14,3).Value = 0
$newSheet.Cells.Item(14,4).Value = 1
$newSheet.Cells.Item(15,1).Value = "Professional Conference"
$newSheet.Cells.Item(15,2).Value = 2
$newSheet.Cells.Item(15,3).Value = 3
$newSheet.Cells.Item(15,4).Value = 0
$newSheet.Cells.Item(16,1).Value = "Purchased List"
$newSheet.Cells.Item(16,2).Value = 13154
$newSheet.Cells.Item(16,3).Value = 16100
$newSheet.Cells.Item(16,4).Value = 10613
$newSheet.Cells.Item(17,1).Value = "Submitted Test Score"
$newSheet.Cells.Item(17,2).Value = 48580
$newSheet.Cells.Item(17,3).Value = 43434
$newSheet.Cells.Item(17,4).Value = 15815
$newSheet.Cells.Item(18,1).Value = "UGrad Campus Visit Check-in Form"
$newSheet.Cells.Item(18,2).Value = 1
$newSheet.Cells.Item(18,3).Value = 0
$newSheet.Cells.Item(18,4).Value = 0
$newSheet.Cells.Item(19,1).Value = "Virtual Tour"
$newSheet.Cells.Item(19,2).Value = 115
$newSheet.Cells.Item(19,3).Value = 812
$newSheet.Cells.Item(19,4).Value = 862
$newSheet.Cells.Item(20,1).Value = "Web Form"
$newSheet.Cells.Item(20,2).Value = 6605
$newSheet.Cells.Item(20,3).Value = 5380
$newSheet.Cells.Item(20,4).Value = 4497

# ---- formatting: reuse the styles already used on the "PRIOR" sheet --
# Row 1 (header) -> bold/bordered/centered style (same as PRIOR!B1)
$ws2.Range("B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Column A (labels) on data rows -> same style as the row-number column
# on the other sheets (PRIOR!A2)
$ws2.Range("A2").Copy()
$newSheet.Range("A2:A20").PasteSpecial(-4122)

# Columns B:D (counts) on data rows -> same fill/border style as PRIOR's
# data row (PRIOR!B2)
$ws2.Range("B2").Copy()
$newSheet.Range("B2:D20").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- column widths (match PRIOR's data-column width for B:D) ---------
$newSheet.Columns.Item(1).ColumnWidth = 25.69
$newSheet.Columns.Item(2).ColumnWidth = $ws2.Columns.Item(2).ColumnWidth
$newSheet.Columns.Item(3).ColumnWidth = $ws2.Columns.Item(2).ColumnWidth
$newSheet.Columns.Item(4).ColumnWidth = $ws2.Columns.Item(2).ColumnWidth

# ---- selection / active view on the new sheet -------------------------
$newSheet.Range("B30").Select()
